# Butman_2007: Updated USDA soil type
#
# Adds a new "pro_usda_soil_order" column to the "profile" sheet (inserted
# between pro_MAP and pro_soil_taxon) and populates the existing entry's
# value. Also adds the corresponding controlled-vocabulary column
# ("pro_usda_soil_order" + the 12 USDA soil order values) to the
# "controlled vocabulary" sheet so the new field has a validation list.

$wb = $excel.ActiveWorkbook

# --- controlled vocabulary sheet -----------------------------------------
$cv = $wb.Worksheets.Item("controlled vocabulary")

# Insert a new column before the existing "pro_soil_taxon_sys" column (E).
$cv.Range("E1").EntireColumn.Insert()

$cv.Range("E1").Value = "pro_usda_soil_order"
$cv.Range("E4").Value = "Alfisols"
$cv.Range("E5").Value = "Andisols"
$cv.Range("E6").Value = "Aridisols"
$cv.Range("E7").Value = "Entisols"
$cv.Range("E8").Value = "Gelisols"
$cv.Range("E9").Value = "Histosols"
$cv.Range("E10").Value = "Inceptisols"
$cv.Range("E11").Value = "Mollisols"
$cv.Range("E12").Value = "Oxisols"
$cv.Range("E13").Value = "Spodosols"
$cv.Range("E14").Value = "Ultisols"
$cv.Range("E15").Value = "Vertisols"

# --- profile sheet -------------------------------------------------------
$profile = $wb.Worksheets.Item("profile")

# Insert a new column before the existing "pro_soil_taxon" column (N).
$profile.Range("N1").EntireColumn.Insert()

$profile.Range("N1").Value = "pro_usda_soil_order"

# Populate the data row (row 4) for this dataset's single profile entry.
$profile.Range("N4").Value = "Inceptisols"
$profile.Range("P4").Value = "Cheshire" + [char]0x2013 + "Holyoke"
$profile.Range("Q4").Value = "USDA"
